$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The attendance sheet used to carry a hard-coded sequence number (1..29) in
# column A. The edit turns that into a "is the typed name the same as the
# roster name" check: column B becomes an editable copy of the roster name
# in column C, and column A becomes a formula that only keeps counting while
# B matches C (otherwise showing "______").

# 1) Duplicate the roster name (column C) into the new column B for every
#    student row (4..32), reusing the existing shared-string entries.
For ($r = 4; $r -le 32; $r++) {
    $ws.Range("B$r").Value = $ws.Range("C$r").Value2
}

# 2) Column A formulas: first row is a plain self-check, second row chains
#    off the first, and the rest (6..31) share one formula referencing the
#    row above - row 32 is left as a literal value, untouched.
$ws.Range("A4").Formula = "=IF(B4=C4,1,""______"")"
$ws.Range("A5").Formula = "=IF(B5=C5,A4+1,""______"")"
$ws.Range("A6:A31").Formula = "=IF(B6=C6,A5+1,""______"")"

# 3) Rows 10 and 13 previously used the special red "missing homework"
#    style (s=8) on column A; restore the normal bordered style (as used by
#    every other row in the column) now that A holds the check formula.
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# 4) Restore the view: frozen panes now show from row 4, and the last
#    selected cell before saving was C31.
[void]$ws.Range("C31").Select()
